# vault backup: 2023-03-27 14:33:29
#
# Update the AutoFilter on column F ("modality") to the new criteria, then
# re-apply the manual row hide/unhide toggles the author made on top of the
# filter result, and leave the selection where the author left it (G3:G72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- AutoFilter --------------------------------------------------------------
# Column D (field 4): drop the "not blank" custom filter entirely.
$ws.Range("A1:G107").AutoFilter(4) | Out-Null

# Column F (field 6, "modality"): now only keep "语音+视频" and
# "语音-视频-文本" (previously "语音", "语音+视频", "语音-视频" plus blanks).
# Re-filtering recomputes which data rows are hidden, so do this before the
# explicit row-visibility overrides below.
$ws.Range("A1:G107").AutoFilter(6, @("语音+视频", "语音-视频-文本"), 7) | Out-Null

# --- Row visibility toggles -------------------------------------------------
# Rows that were previously hidden by the filter and are now shown again.
$ws.Rows.Item(3).Hidden  = $false
$ws.Rows.Item(22).Hidden = $false
$ws.Rows.Item(51).Hidden = $false
$ws.Rows.Item(64).Hidden = $false
$ws.Rows.Item(72).Hidden = $false

# Rows that were previously visible and are now manually hidden.
$ws.Rows.Item(85).Hidden  = $true
$ws.Rows.Item(96).Hidden  = $true
$ws.Rows.Item(97).Hidden  = $true
$ws.Rows.Item(106).Hidden = $true
$ws.Rows.Item(107).Hidden = $true

# --- Selection ---------------------------------------------------------------
$ws.Range("G3:G72").Select() | Out-Null
